$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (source cell to copy number-format/style from, value to write)
# Column K holds the "style donor" cells that already carry the right xf
# (same visual style the new column M should use, mirroring column L's pattern).
$rows = @(
    @{ Row = 4;  Src = "L4";  Val = 2022 },
    @{ Row = 5;  Src = "K5";  Val = 24.6 },
    @{ Row = 6;  Src = "K6";  Val = 40.700000000000003 },
    @{ Row = 7;  Src = "K7";  Val = 20.7 },
    @{ Row = 8;  Src = "K8";  Val = 26.6 },
    @{ Row = 9;  Src = "K9";  Val = 44.5 },
    @{ Row = 10; Src = "K10"; Val = 21.9 },
    @{ Row = 11; Src = "K11"; Val = 21.9 },
    @{ Row = 12; Src = "K12"; Val = 35.299999999999997 },
    @{ Row = 13; Src = "K13"; Val = 17.600000000000001 },
    @{ Row = 15; Src = "K15"; Val = 44.9 },
    @{ Row = 16; Src = "K16"; Val = 21.5 },
    @{ Row = 17; Src = "K17"; Val = 36.200000000000003 },
    @{ Row = 18; Src = "K18"; Val = 53.1 },
    @{ Row = 19; Src = "K19"; Val = 33.4 },
    @{ Row = 20; Src = "K20"; Val = 20.2 },
    @{ Row = 21; Src = "K21"; Val = 15.4 },
    @{ Row = 22; Src = "K22"; Val = 20.5 },
    @{ Row = 23; Src = "K23"; Val = 27.1 },
    @{ Row = 24; Src = "K24"; Val = 36.1 },
    @{ Row = 25; Src = "K25"; Val = 25.2 },
    @{ Row = 26; Src = "K26"; Val = 24.2 },
    @{ Row = 27; Src = "K27"; Val = 46.5 },
    @{ Row = 28; Src = "K28"; Val = 20.3 },
    @{ Row = 29; Src = "K29"; Val = 40.5 },
    @{ Row = 30; Src = "L30"; Val = 44.5 }
)

foreach ($item in $rows) {
    $target = "M" + $item.Row
    $ws.Range($item.Src).Copy()
    $ws.Range($target).PasteSpecial(-4122)
    $ws.Range($target).Value = $item.Val
}

# Row 14 needs a brand-new style: same as the bold "K-column" style (fontId 2,
# no border/fill) but with the 0.0 number format applied on top of it -- this
# mints a new cellXfs entry rather than reusing one.
$ws.Range("K14").Copy()
$ws.Range("M14").PasteSpecial(-4122)
$ws.Range("M14").NumberFormat = "0.0"
$ws.Range("M14").Value = 28

# Update the active selection to match the saved view state.
[void]$ws.Range("N7").Select()

Write-Output "done"
